$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "303.04"
Set-TextValue "E2" "-0.39%"
Set-TextValue "D3" "37.49"
Set-TextValue "E3" "6.47%"
Set-TextValue "E4" "-3.29%"
Set-TextValue "D5" "0.07872"
Set-TextValue "E5" "0.36%"
Set-TextValue "D6" "2.213"
Set-TextValue "E6" "-3.68%"
Set-TextValue "D7" "8.004"
Set-TextValue "E7" "-0.65%"
Set-TextValue "D8" "4.029"
Set-TextValue "E8" "1.05%"
Set-TextValue "D9" "0.9162"
Set-TextValue "D10" "0.09558"
Set-TextValue "E10" "-4.71%"
Set-TextValue "D11" "0.1875"
Set-TextValue "E11" "2.11%"
Set-TextValue "D12" "0.08569"
Set-TextValue "E12" "0.17%"
Set-TextValue "D13" "0.03556"
Set-TextValue "E13" "5.60%"
Set-TextValue "D14" "0.09955"
Set-TextValue "E14" "0.46%"
Set-TextValue "D15" "0.001469"
Set-TextValue "E15" "-0.64%"
Set-TextValue "D16" "0.005702"
Set-TextValue "E16" "-0.97%"
Set-TextValue "E17" "-0.65%"
Set-TextValue "D18" "2.250"
Set-TextValue "E18" "6.59%"
Set-TextValue "D20" "0.1317"
Set-TextValue "E20" "-0.66%"
Set-TextValue "D21" "4.765"
Set-TextValue "E21" "4.34%"
Set-TextValue "D22" "0.2223"
Set-TextValue "E22" "-6.56%"
Set-TextValue "D23" "0.04583"
Set-TextValue "E23" "-1.39%"
Set-TextValue "D24" "0.001231"
Set-TextValue "E24" "0.71%"
Set-TextValue "D25" "0.004785"
Set-TextValue "E25" "7.43%"
Set-TextValue "D26" "0.0001400"
Set-TextValue "E26" "7.93%"
Set-TextValue "E27" "40.05%"
Set-TextValue "D39" "0.01820"
Set-TextValue "E39" "4.07%"
Set-TextValue "D40" "0.04727"
Set-TextValue "E40" "-0.38%"
Set-TextValue "D41" "0.008111"
Set-TextValue "E41" "5.55%"
Set-TextValue "D42" "0.1393"
Set-TextValue "E42" "-1.40%"
Set-TextValue "D43" "0.007559"
Set-TextValue "E43" "7.37%"
Set-TextValue "D44" "0.002210"
Set-TextValue "E44" "-3.70%"
Set-TextValue "D45" "0.01044"
Set-TextValue "E45" "5.64%"
Set-TextValue "D46" "0.00006159"
Set-TextValue "E46" "2.87%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "0.23%"
Set-TextValue "D48" "0.0005801"
Set-TextValue "E48" "0.01%"
Set-TextValue "D49" "6.678"
Set-TextValue "E49" "16.27%"
Set-TextValue "E50" "0.24%"
Set-TextValue "D51" "0.00002100"
Set-TextValue "E51" "0.23%"
